# Update close-contact residue force data rows (within 6A) in the
# permeation_frames_forces sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> column letter -> new value
$changes = @{
    8  = @{
        B = 4415
        C = '[-1.2431654334068298, -0.5947210304439068, -11.592971801757812]'
        D = 11.6745941511995
        E = 11.5522255221292
        F = 0.9895183826105229
        G = 1.378097746486022
        H = -11.59297180175781
        I = '[-1.1576690673828125, 0.13907241821289062, -4.969749450683594]'
    }
    9  = @{
        B = 4994
        C = '[-1.2688956409692764, -0.39530257135629654, -15.477433919906616]'
        D = 15.53439155601749
        E = 15.31075270340408
        F = 0.9856036297394098
        G = 1.329044871549389
        H = -15.47743391990662
        I = '[0.37653350830078125, -0.006805419921875, -4.346221923828125]'
    }
    11 = @{
        B = 5399
        C = '[-2.5655597448349, 2.268065929412842, -10.400489449501038]'
        D = 10.94972148749287
        E = 9.907631183671697
        F = 0.9048295150692662
        G = 3.424356854137818
        H = -10.40048944950104
        I = '[0.17792129516601562, -0.7498664855957031, -5.464141845703125]'
    }
    16 = @{
        B = 6016
        C = '[-4.774939412251115, 1.9371942728757858, -10.307047605514526]'
        D = 11.52336749322326
        E = 9.361208832082269
        F = 0.8123674644228325
        G = 5.152937806875973
        H = -10.30704760551453
        I = '[0.917236328125, -0.059780120849609375, -5.506996154785156]'
    }
    19 = @{
        B = 6426
        C = '[-1.5980022549629211, -2.4518961906433105, -5.547749876976013]'
        D = 6.272394665041823
        E = 2.890896342806211
        F = 0.4608919714376639
        G = 2.926671511556731
        H = -5.547749876976013
        I = '[1.1834907531738281, 0.7636146545410156, -1.9186553955078125]'
    }
    20 = @{
        B = 6488
        C = '[-0.4124993681907654, -1.4086133688688278, -8.854763507843018]'
        D = 8.975588233176818
        E = 5.568802080930866
        F = 0.6204386761356414
        G = 1.467769515868813
        H = -8.854763507843018
        I = '[-0.6596145629882812, 1.8584747314453125, -2.0570907592773438]'
    }
    21 = @{
        B = 6561
        C = '[-2.351348042488098, 2.377216547727585, -11.267191410064697]'
        D = 11.75285490435146
        E = 10.09789668287035
        F = 0.8591867052771691
        G = 3.343650121005855
        H = -11.2671914100647
        I = '[-0.17771148681640625, -1.4416885375976562, -4.5340423583984375]'
    }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
